$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: E9 "D" -> "W", F9 "W" -> "D" (values swapped)
$ws.Range("E9").Value = "W"
$ws.Range("F9").Value = "D"

# Row 10: D10 "W" -> "_", E10 "_" -> "W" (values swapped)
$ws.Range("D10").Value = "_"
$ws.Range("E10").Value = "W"

# Row 11: C11 "_" -> "D", D11 "W" -> "_", E11 "_" -> "W"
$ws.Range("C11").Value = "D"
$ws.Range("D11").Value = "_"
$ws.Range("E11").Value = "W"

# Update the view selection to P9 (also drops the stale topLeftCell scroll anchor)
$ws.Range("P9").Select()
